$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "3.0-Leche (litros),2.0-Harina  (kg),1.0-Vainilla (ml),1.0-Huevos (unidad),"
$ws.Range("C3").Value = "2.0-Harina  (kg),5.0-Huevos (unidad),"
$ws.Range("C4").Value = "5.0-Harina  (kg),1.0-Vainilla (ml),2.0-Huevos (unidad),"
$ws.Range("C5").Value = "5.0-Harina  (kg),5.0-Huevos (unidad),"
$ws.Range("C6").Value = "5.0-Crema (litros),4.0-Harina  (kg),5.0-Huevos (unidad),2.0-Limon (unidad),"
$ws.Range("C7").Value = "0.2-Leche (litros),0.3-Harina  (kg),0.1-Vainilla (ml),2.0-Huevos (unidad),"
